# Insert a new data row at row 49 (pushing existing rows 49..161 down to 50..162),
# and populate the newly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(49).Insert()

$ws.Cells.Item(49, 1).Value = 5
$ws.Cells.Item(49, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(49, 3).Value = 'Maule'
$ws.Cells.Item(49, 4).Value = 45281
$ws.Cells.Item(49, 5).Value = 7
$ws.Cells.Item(49, 6).Value = 'Fruta'
$ws.Cells.Item(49, 7).Value = 100101
$ws.Cells.Item(49, 8).Value = 'Berries'
$ws.Cells.Item(49, 9).Value = 100101001
$ws.Cells.Item(49, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(49, 11).Value = 'Sin especificar'
$ws.Cells.Item(49, 12).Value = 'Primera'
$ws.Cells.Item(49, 13).Value = 130
$ws.Cells.Item(49, 14).Value = 4000
$ws.Cells.Item(49, 15).Value = 4000
$ws.Cells.Item(49, 16).Value = 4000
$ws.Cells.Item(49, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(49, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(49, 19).Value = 2000
$ws.Cells.Item(49, 20).Value = 2
